$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C59").Value = "Merry Christmas & Happy New Year Itsuki! We loving having you at Super Epion, you're always so happy and well behaved in class, and your English is improving so much already! Good job!"
$ws.Range("C61").Value = "Merry Christmas & Happy New Year Hinata! We loving having you at Super Epion, you're always so happy and well behaved in class, and your English is improving so much already! Good job!"
$ws.Range("C69").Value = "Merry Christmas & Happy New Year Miyu! We loving having you at Super Epion, you're always so happy and well behaved in class, and your English is improving so much already! Good job!"

$ws.Range("C55").Select()
$ws.Application.ActiveWindow.ScrollRow = 30
